# Generate Report for Handoff
#
# Refreshes the localization-status report for a new handoff run:
#   - the source markdown file got a new GUID-based name
#     (a41a0594-4834-42c3-8ca6-855dc2410c7a -> a033269e-21b4-4a37-b892-79af958f5f0a)
#   - the generated xliff files got a new content hash
#     (94a1e704d5d8ff888fad466f081da08c713bf463 -> c55a0630db6d8c71a6eb18ee157528609ee2e6db)
#   - timestamps advance to reflect the new run

$wb = $excel.ActiveWorkbook

$oldGuid = "a41a0594-4834-42c3-8ca6-855dc2410c7a"
$newGuid = "a033269e-21b4-4a37-b892-79af958f5f0a"
$oldHash = "94a1e704d5d8ff888fad466f081da08c713bf463"
$newHash = "c55a0630db6d8c71a6eb18ee157528609ee2e6db"

$newFileName = "$newGuid.md"
$newPathAndName = "e2e\$newGuid.md"

# Hyperlinks keep pointing at the same external target (github blob URL);
# only the cell's displayed text changes. This host re-creates (rather than
# mutates) a Hyperlink when any of its properties are written, so delete +
# re-add to avoid leaving a stale duplicate hyperlink behind.
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c89bcd93c1e9557313cd9785264e2726bb92749b/e2e/$oldGuid.md"

function Update-HyperlinkDisplay($range, $displayText) {
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $hyperlinkTarget, $null, $null, $displayText) | Out-Null
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
Update-HyperlinkDisplay $wsOverview.Range("B2") $newPathAndName
$wsOverview.Range("G2").Value = "2016-09-05 09:19:39"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newFileName
Update-HyperlinkDisplay $wsZhCn.Range("A2") $newFileName
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-05 09:19:33"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newFileName
Update-HyperlinkDisplay $wsDeDe.Range("A2") $newFileName
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-05 09:19:39"
